$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (row 878 through 920) appended to the close/train scaler data,
# containing the rerun DLinear output on the updated dataset.
$newData = @(
    @(876, 0.88049),
    @(877, 0.8792),
    @(878, 0.88195),
    @(879, 0.8762799999999999),
    @(880, 0.87093),
    @(881, 0.84612),
    @(882, 0.85017),
    @(883, 0.85017),
    @(884, 0.85017),
    @(885, 0.85439),
    @(886, 0.84612),
    @(887, 0.84206),
    @(888, 0.82536),
    @(889, 0.83363),
    @(890, 0.83785),
    @(891, 0.8226),
    @(892, 0.83104),
    @(893, 0.82406),
    @(894, 0.81158),
    @(895, 0.79796),
    @(896, 0.76342),
    @(897, 0.74964),
    @(898, 0.79796),
    @(899, 0.74964),
    @(900, 0.75661),
    @(901, 0.76893),
    @(902, 0.74282),
    @(903, 0.70829),
    @(904, 0.70829),
    @(905, 0.66694),
    @(906, 0.68623),
    @(907, 0.6945),
    @(908, 0.70829),
    @(909, 0.70569),
    @(910, 0.72077),
    @(911, 0.72904),
    @(912, 0.76066),
    @(913, 0.70699),
    @(914, 0.6684),
    @(915, 0.69726),
    @(916, 0.70569),
    @(917, 0.73326),
    @(918, 0.72207)
)

$startRow = 878
$lastExistingRow = 877

# Propagate the existing data-row formatting (border + centered style) from
# the last populated row down across the newly added rows in column A.
$endRow = $startRow + $newData.Count - 1
$ws.Range("A" + $lastExistingRow).Copy($ws.Range("A" + $startRow + ":A" + $endRow))

$r = $startRow
foreach ($pair in $newData) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}
